$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.579.66"
$ws.Range("E2").Value = "  -3.87%  "
$ws.Range("D3").Value = "2.973.17"
$ws.Range("E3").Value = "  -5.14%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "542.01"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "151.92"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.01%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  +0.40%  "
$ws.Range("D9").Value = "2.981.07"
$ws.Range("E9").Value = "  -5.27%  "
$ws.Range("E10").Value = "  -2.18%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.13"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -6.76%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.369"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.74%  "
$ws.Range("D13").Value = "3.493.79"
$ws.Range("E13").Value = "  -5.23%  "
$ws.Range("E14").Value = "  -2.30%  "
$ws.Range("D15").Value = "61.663.22"
$ws.Range("E15").Value = "  -3.85%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "23.73"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.18%  "
$ws.Range("D17").Value = "2.976.25"
$ws.Range("E17").Value = "  -5.36%  "
$ws.Range("E18").Value = "  -4.08%  "
$ws.Range("E19").Value = "  -1.18%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.02"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.95%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "381.03"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.67"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -6.21%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.65"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.46%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.69"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.33%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.472"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.38%  "
$ws.Range("D27").Value = "3.102.40"
$ws.Range("E27").Value = "  -5.22%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.189"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.98%  "
$ws.Range("D29").Value = "0.0₃0946"
$ws.Range("E29").Value = "  -5.60%  "
$ws.Range("E30").Value = "  +1.06%  "
$ws.Range("E31").Value = "  -5.63%  "
$ws.Range("E32").Value = "  +0.02%  "
$ws.Range("E33").Value = "  -4.61%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "20.53"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.72%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "160.77"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.87%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.67"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.79%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.92"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.19%  "
$ws.Range("E38").Value = "  -2.78%  "
$ws.Range("E39").Value = "  -5.19%  "
$ws.Range("E40").Value = "  -6.30%  "
$ws.Range("E41").Value = "  -2.94%  "
$ws.Range("B42").Value = "OKB"
$ws.Range("C42").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "37.44"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.21%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "2.414.24"
$ws.Range("E43").Value = "  -9.37%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "22.16"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.14%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.669"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.86%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0594"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.08%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.13"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.24%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0247"
$ws.Range("D48").Style = "Normal"
$ws.Range("E49").Value = "  +0.00%  "
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0954"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.02%  "
$ws.Range("B51").Value = "Bittensor"
$ws.Range("C51").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "267.49"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -6.60%  "
